{"js": "// Translate the English subtitle text to Italian (video_subtitles/translation/ita).\n// Each entry is the unique, exact text of one paragraph's run in before.docx.\nconst replacements = [\n  [\"Format has been corrected not the timing\", \"Il formato \u00e8 stato corretto, non la tempistica\"],\n  [\"I added 25 seconds to each timing to correct for the intro song -john argentino\", \"Ho aggiunto 25 secondi a ogni periodi per correggere per la canzone di introduzione -John Argentino\"],\n  [\"The airport problem - subtitles:\", \"Il problema dell'aeroporto - sottotitoli:\"],\n  [\"The administrations of three\", \"Le amministrazioni di tre\"],\n  [\"neighboring cities: A, B and C decided\", \"citt\u00e0 confinanti: A, B e C hanno deciso\"],\n  [\"to build an airport dividing the costs of\", \"di costruire un aeroporto che divida i costi di\"],\n  [\"implementation. The condition on the\", \"implementazione. La condizione sulla\"],\n  [\"choice of the most suitable place is\", \"scelta del posto pi\u00f9 adatto \u00e8\"],\n  [\"that the sum of the distances from each\", \"che la somma delle distanze da ogni\"],\n  [\"city to the airport is as small as\", \"citt\u00e0 all'aeroporto sia la minore\"],\n  [\"possible. The team of experts in charge\", \"possibile. Il team di esperti in carica\"],\n  [\"of the work has created a model to get\", \"del lavoro ha creato un modello per ottenere\"],\n  [\"a preliminary idea of where to place the\", \"un'idea preliminare di dove posizionare la\"],\n  [\"structure. At their disposal there are\", \"struttura. A loro disposizione ci sono\"],\n  [\"some snails a big metal ring and a long\", \"dei chiodi, un grande anello di metallo e una lunga\"],\n  [\"string.\", \"corda.\"],\n  [\"Explain how the team can manage to use\", \"Spiega come il team possa usare i\"],\n  [\"the materials to tell approximately the\", \"materiali per dire approssimativamente la\"],\n  [\"ideal location of the airport. Imagine\", \"posizione ideale dell'aeroporto. Immagina\"],\n  [\"that the cities are placed at the\", \"che le citt\u00e0 si trovino ai\"],\n  [\"vertices of a triangle which is\", \"vertici di un triangolo che sia\"],\n  [\"obviously reproduced in scale as\", \"ovviamente riprodotto in scala come\"],\n  [\"shown in figure. This is one possible\", \"mostrato in figura. This is one possible\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Translate the English subtitle text to Italian (video_subtitles/translation/ita).\n# Each \"Find\" string is the unique, exact text of one paragraph's run in before.docx.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Format has been corrected not the timing\", \"Il formato \u00e8 stato corretto, non la tempistica\"),\n    @(\"I added 25 seconds to each timing to correct for the intro song -john argentino\", \"Ho aggiunto 25 secondi a ogni periodi per correggere per la canzone di introduzione -John Argentino\"),\n    @(\"The airport problem - subtitles:\", \"Il problema dell'aeroporto - sottotitoli:\"),\n    @(\"The administrations of three\", \"Le amministrazioni di tre\"),\n    @(\"neighboring cities: A, B and C decided\", \"citt\u00e0 confinanti: A, B e C hanno deciso\"),\n    @(\"to build an airport dividing the costs of\", \"di costruire un aeroporto che divida i costi di\"),\n    @(\"implementation. The condition on the\", \"implementazione. La condizione sulla\"),\n    @(\"choice of the most suitable place is\", \"scelta del posto pi\u00f9 adatto \u00e8\"),\n    @(\"that the sum of the distances from each\", \"che la somma delle distanze da ogni\"),\n    @(\"city to the airport is as small as\", \"citt\u00e0 all'aeroporto sia la minore\"),\n    @(\"possible. The team of experts in charge\", \"possibile. Il team di esperti in carica\"),\n    @(\"of the work has created a model to get\", \"del lavoro ha creato un modello per ottenere\"),\n    @(\"a preliminary idea of where to place the\", \"un'idea preliminare di dove posizionare la\"),\n    @(\"structure. At their disposal there are\", \"struttura. A loro disposizione ci sono\"),\n    @(\"some snails a big metal ring and a long\", \"dei chiodi, un grande anello di metallo e una lunga\"),\n    @(\"string.\", \"corda.\"),\n    @(\"Explain how the team can manage to use\", \"Spiega come il team possa usare i\"),\n    @(\"the materials to tell approximately the\", \"materiali per dire approssimativamente la\"),\n    @(\"ideal location of the airport. Imagine\", \"posizione ideale dell'aeroporto. Immagina\"),\n    @(\"that the cities are placed at the\", \"che le citt\u00e0 si trovino ai\"),\n    @(\"vertices of a triangle which is\", \"vertici di un triangolo che sia\"),\n    @(\"obviously reproduced in scale as\", \"ovviamente riprodotto in scala come\"),\n    @(\"shown in figure. This is one possible\", \"mostrato in figura. This is one possible\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $findText\n    $rng.Find.Replacement.Text = $replaceText\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n\n    $found = $rng.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: not found -> $findText\"\n    }\n}\n"}
